$d = $word.ActiveDocument

# Keep straight apostrophes as typed (defensive: avoid Word's "smart quotes"
# autocorrect mangling text we insert/replace below).
$word.Options.AutoFormatAsYouTypeReplaceQuotes = $false
$word.Options.AutoFormatReplaceQuotes = $false

$pkgHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + "`n" + `
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
    '<pkg:xmlData>' + `
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
    '<w:body>'
$pkgFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# 1. Remove the "Meta description" paragraph (2nd paragraph of the document,
#    right after the Heading1 title) entirely, including its paragraph mark.
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Delete()

# 2. Insert a new bold paragraph, with the same title text as the Heading1,
#    right before the final "Prompt: ..." paragraph (i.e. just before
#    the section break / end of document).
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$lastPara.Range.InsertParagraphBefore() | Out-Null

$newPara = $d.Paragraphs.Item($count)
$titleXml = $pkgHeader + `
    '<w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Finn''s Golden Tavern for Free - Innovative Spiral Grid Gameplay</w:t></w:r></w:p>' + `
    $pkgFooter
$newPara.Range.InsertXML($titleXml) | Out-Null

# 3. Replace the content of the final "Prompt: ..." paragraph's italic run
#    with the meta-description text (minus the "Meta description: " label),
#    leaving its leading empty run and italic formatting intact. We target
#    the paragraph content excluding the trailing paragraph mark, and the
#    replacement XML excludes its own leading empty run (the paragraph's
#    existing one is preserved untouched) so we don't end up with a
#    duplicate empty run.
$count2 = $d.Paragraphs.Count
$promptPara = $d.Paragraphs.Item($count2)
$fullRange = $promptPara.Range
$contentRange = $d.Range($fullRange.Start, $fullRange.End - 1)
$descXml = $pkgHeader + `
    '<w:p><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve">Read our review of Finn''s Golden Tavern and play it for free. Innovative spiral grid gameplay with excellent graphics and original features. </w:t></w:r></w:p>' + `
    $pkgFooter
$contentRange.InsertXML($descXml) | Out-Null
